# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the country names for rows 12 and 13 (Rusia overtakes China in ranking)
$ws.Range("A12").Value = "Rusia"
$ws.Range("A13").Value = "China"

# Row 12 (Rusia) - updated statistics
$ws.Range("B12").Value = 87147
$ws.Range("C12").Value = 6198
$ws.Range("D12").Value = 7346
$ws.Range("E12").Value = 79007
$ws.Range("F12").Value = 2300
$ws.Range("G12").Value = 47
$ws.Range("H12").Value = 794

# Row 13 (China) - takes on the previous China statistics (row 12's old values)
$ws.Range("B13").Value = 82830
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 77474
$ws.Range("E13").Value = 723
$ws.Range("F13").Value = 52
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 4633

# Row 33 (Polonia) - updated statistics
$ws.Range("B33").Value = 11761
$ws.Range("C33").Value = 144
$ws.Range("E33").Value = 8756
$ws.Range("G33").Value = 4
$ws.Range("H33").Value = 539

# Row 34 (Rumania) - updated statistics
$ws.Range("E34").Value = 7354
$ws.Range("G34").Value = 9
$ws.Range("H34").Value = 628

# Row 41 (Dinamarca) - updated statistics
$ws.Range("B41").Value = 8698
$ws.Range("C41").Value = 123
$ws.Range("E41").Value = 2471

# Row 43 (Filipinas) - updated statistics
$ws.Range("B43").Value = 7777
$ws.Range("C43").Value = 198
$ws.Range("D43").Value = 932
$ws.Range("E43").Value = 6334
$ws.Range("F43").Value = 31
$ws.Range("G43").Value = 10
$ws.Range("H43").Value = 511

# Row 62 (Kazajistan) - updated statistics
$ws.Range("B62").Value = 2791
$ws.Range("C62").Value = 74
$ws.Range("E62").Value = 2084
